$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "Launch Cash Builder Page" marked Completed, with a completed date ---
$ws.Range("B6").Interior.Color = 5296274   # green fill (RGB 92,208,80) - same as "Completed" key swatch
$ws.Range("C6").Value = 43867              # 06/02/2020 (Excel serial date)
$ws.Range("C6").NumberFormat = "mm-dd-yy"  # built-in short-date format (numFmtId 14)

# --- Row 8: "Start Timer" marked Completed, with a completed date + note ---
$ws.Range("B8").Interior.Color = 5296274   # green fill - same as "Completed" key swatch

# Copy C6's date formatting (value + number format) onto C8 so both cells share
# a single reused style entry, exactly like the recorded styles would.
$ws.Range("C6").Copy()
$ws.Range("C8").PasteSpecial(-4122)        # xlPasteFormats
$ws.Range("C8").Value = 43867              # 06/02/2020 (Excel serial date)

$ws.Range("D8").Value = "Only basic implementation"

# New column D needs to fit the note text added above
$ws.Columns.Item(4).AutoFit()

# Update the active selection left behind after editing
$ws.Range("C11").Select()
